# Add a "Channel" column in front of the existing data (channel column in write module added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A, shifting N / Inversion / Numero de combinaciones / Cobertura one column right
$ws.Columns("A").Insert()

# Header for the new column
$ws.Range("A1").Value = "Channel"

# Channel names per row
$ws.Range("A2").Value = "Televisión"
$ws.Range("A3").Value = "Digital Video"
$ws.Range("A4").Value = "Cine"
$ws.Range("A5").Value = "BVOD"

# Updated Inversión figures (now column C after the insert)
$ws.Range("C2").Value = 443750
$ws.Range("C3").Value = 525000
$ws.Range("C4").Value = 231250
$ws.Range("C5").Value = 0

# Updated Cobertura figure recalculated for the new Inversión value (now column E)
$ws.Range("E2").Value = 0.8274212228066136
